$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 501716.16
$ws.Range("I40").Value = 1355.5555
$ws.Range("J40").Value = 911102.0600000001
$ws.Range("K40").Value = 1355.5555
$ws.Range("L40").Value = 911102.0600000001
$ws.Range("M40").Value = -1180.5555
$ws.Range("N40").Value = -911452.0600000001

# Row 58
$ws.Range("H58").Value = 287
$ws.Range("I58").Value = 108.75
$ws.Range("J58").Value = 1000
$ws.Range("K58").Value = 326.25
$ws.Range("L58").Value = 3000
$ws.Range("M58").Value = -176.25
$ws.Range("N58").Value = -3300

# Row 62
$ws.Range("H62").Value = 20847834
$ws.Range("I62").Value = 27791068
$ws.Range("K62").Value = 27791068
$ws.Range("M62").Value = -27790444

# Row 65
$ws.Range("H65").Value = 20847834
$ws.Range("I65").Value = 27791068
$ws.Range("K65").Value = 138955340
$ws.Range("M65").Value = -138952220

# Row 98
$ws.Range("H98").Value = 65790670
$ws.Range("I98").Value = 78126136
$ws.Range("J98").Value = 1500
$ws.Range("K98").Value = 78126136
$ws.Range("L98").Value = 1500
$ws.Range("M98").Value = -78124638
$ws.Range("N98").Value = -4496

# Row 113
$ws.Range("H113").Value = 2887.3948
$ws.Range("I113").Value = 4363
$ws.Range("J113").Value = 2663.818
$ws.Range("K113").Value = 4363
$ws.Range("L113").Value = 2663.818
$ws.Range("M113").Value = -1109
$ws.Range("N113").Value = -9171.817999999999

# Row 122
$ws.Range("H122").Value = 65790670
$ws.Range("I122").Value = 78126136
$ws.Range("J122").Value = 1500
$ws.Range("K122").Value = 234378408
$ws.Range("L122").Value = 4500
$ws.Range("M122").Value = -234375958
$ws.Range("N122").Value = -9400

# Row 132
$ws.Range("H132").Value = 22739818
$ws.Range("I132").Value = 23822638
$ws.Range("J132").Value = 600
$ws.Range("K132").Value = 71467914
$ws.Range("L132").Value = 1800
$ws.Range("M132").Value = -71465384
$ws.Range("N132").Value = -6860

# Row 138
$ws.Range("H138").Value = 5197.857
$ws.Range("I138").Value = 1738.7368
$ws.Range("J138").Value = 9305.5625
$ws.Range("K138").Value = 5216.2104
$ws.Range("L138").Value = 27916.6875
$ws.Range("M138").Value = -76.21039999999994
$ws.Range("N138").Value = -38196.6875

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 10109986
$ws.Range("I32").Value = 8060.923
$ws.Range("K32").Value = 8060.923
$ws.Range("M32").Value = -7773.923

# Row 61
$ws.Range("H61").Value = 8335866
$ws.Range("I61").Value = 12502020
$ws.Range("J61").Value = 3558.8
$ws.Range("K61").Value = 12502020
$ws.Range("L61").Value = 3558.8
$ws.Range("M61").Value = -12501808
$ws.Range("N61").Value = -3982.8

# Row 122
$ws.Range("H122").Value = 2113.9
$ws.Range("I122").Value = 1892.375
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 5677.125
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -3227.125
$ws.Range("N122").Value = -13900

# Row 123
$ws.Range("H123").Value = 28000
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 28000
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 28000
$ws.Range("M123").ClearContents()
$ws.Range("N123").Value = -37800

# Row 132
$ws.Range("H132").Value = 2264556.8
$ws.Range("I132").Value = 2117.4
$ws.Range("J132").Value = 9806021
$ws.Range("K132").Value = 6352.200000000001
$ws.Range("L132").Value = 29418063
$ws.Range("M132").Value = -3822.200000000001
$ws.Range("N132").Value = -29423123

# Row 136
$ws.Range("H136").Value = 8335866
$ws.Range("I136").Value = 12502020
$ws.Range("J136").Value = 3558.8
$ws.Range("K136").Value = 37506060
$ws.Range("L136").Value = 10676.4
$ws.Range("M136").Value = -37503510
$ws.Range("N136").Value = -15776.4

$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 1058.8636
$ws.Range("I94").Value = 894.41174
$ws.Range("J94").Value = 1618
$ws.Range("K94").Value = 894.41174
$ws.Range("L94").Value = 1618
$ws.Range("M94").Value = -443.41174
$ws.Range("N94").Value = -2520

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1455.875
$ws.Range("I31").Value = 880.3684
$ws.Range("J31").Value = 3642.8
$ws.Range("K31").Value = 880.3684
$ws.Range("L31").Value = 3642.8
$ws.Range("M31").Value = -585.3684
$ws.Range("N31").Value = -4232.8

# Row 34
$ws.Range("H34").Value = 1455.875
$ws.Range("I34").Value = 880.3684
$ws.Range("J34").Value = 3642.8
$ws.Range("K34").Value = 880.3684
$ws.Range("L34").Value = 3642.8
$ws.Range("M34").Value = -678.3684
$ws.Range("N34").Value = -4046.8

# Row 51
$ws.Range("H51").Value = 43320
$ws.Range("J51").Value = 43320
$ws.Range("L51").Value = 43320
$ws.Range("N51").Value = -44792

# Row 58
$ws.Range("H58").Value = 33334370
$ws.Range("I58").Value = 50000900
$ws.Range("J58").Value = 1310
$ws.Range("K58").Value = 50000900
$ws.Range("L58").Value = 1310
$ws.Range("M58").Value = -50000697
$ws.Range("N58").Value = -1716

# Row 59
$ws.Range("H59").Value = 38756
$ws.Range("J59").Value = 38756
$ws.Range("L59").Value = 38756
$ws.Range("N59").Value = -41046

# Row 61
$ws.Range("H61").Value = 43320
$ws.Range("J61").Value = 43320
$ws.Range("L61").Value = 43320
$ws.Range("N61").Value = -44016

# Row 74
$ws.Range("H74").Value = 20555.334
$ws.Range("J74").Value = 18333
$ws.Range("L74").Value = 18333
$ws.Range("N74").Value = -20081

# Row 77
$ws.Range("H77").Value = 20555.334
$ws.Range("J77").Value = 18333
$ws.Range("L77").Value = 54999
$ws.Range("N77").Value = -63735

# Row 92
$ws.Range("H92").Value = 20300.5
$ws.Range("J92").Value = 20300.5
$ws.Range("L92").Value = 20300.5
$ws.Range("N92").Value = -25292.5

# Row 122
$ws.Range("H122").Value = 13159266
$ws.Range("I122").Value = 19232276
$ws.Range("J122").Value = 1078
$ws.Range("K122").Value = 57696828
$ws.Range("L122").Value = 3234
$ws.Range("M122").Value = -57694378
$ws.Range("N122").Value = -8134

# Row 132
$ws.Range("H132").Value = 12822381
$ws.Range("I132").Value = 1551
$ws.Range("J132").Value = 47621776
$ws.Range("K132").Value = 4653
$ws.Range("L132").Value = 142865328
$ws.Range("M132").Value = -2123
$ws.Range("N132").Value = -142870388

# Row 134
$ws.Range("H134").Value = 1327.3928
$ws.Range("I134").Value = 1069.6111
$ws.Range("J134").Value = 1791.4
$ws.Range("K134").Value = 3208.8333
$ws.Range("L134").Value = 5374.200000000001
$ws.Range("M134").Value = -673.8333000000002
$ws.Range("N134").Value = -10444.2

# Row 136
$ws.Range("H136").Value = 33334370
$ws.Range("I136").Value = 50000900
$ws.Range("J136").Value = 1310
$ws.Range("K136").Value = 150002700
$ws.Range("L136").Value = 3930
$ws.Range("M136").Value = -150000150
$ws.Range("N136").Value = -9030

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 38892948
$ws.Range("J5").Value = 7233.625
$ws.Range("L5").Value = 21700.875
$ws.Range("N5").Value = -21924.875

# Row 122
$ws.Range("H122").Value = 10643248
$ws.Range("I122").Value = 62500460
$ws.Range("J122").Value = 5871.641
$ws.Range("K122").Value = 562504140
$ws.Range("L122").Value = 52844.769
$ws.Range("M122").Value = -562501690
$ws.Range("N122").Value = -57744.769

# Row 132
$ws.Range("H132").Value = 8467.666999999999
$ws.Range("J132").Value = 8965.357
$ws.Range("L132").Value = 80688.213
$ws.Range("N132").Value = -85748.213

# Row 135
$ws.Range("H135").Value = 38892948
$ws.Range("J135").Value = 7233.625
$ws.Range("L135").Value = 65102.625
$ws.Range("N135").Value = -70172.625

$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 2024
$ws.Range("I102").Value = 1593.6
$ws.Range("J102").Value = 3100
$ws.Range("K102").Value = 1593.6
$ws.Range("L102").Value = 3100
$ws.Range("M102").Value = 28.40000000000009
$ws.Range("N102").Value = -6344

# Row 122
$ws.Range("H122").Value = 41684450
$ws.Range("I122").Value = 41684450
$ws.Range("K122").Value = 125053350
$ws.Range("M122").Value = -125050900

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 1267.6471
$ws.Range("I22").Value = 2133.3333
$ws.Range("J22").Value = 795.4545000000001
$ws.Range("K22").Value = 2133.3333
$ws.Range("L22").Value = 795.4545000000001
$ws.Range("M22").Value = -1838.3333
$ws.Range("N22").Value = -1385.4545

# Row 27
$ws.Range("H27").Value = 1267.6471
$ws.Range("I27").Value = 2133.3333
$ws.Range("J27").Value = 795.4545000000001
$ws.Range("K27").Value = 2133.3333
$ws.Range("L27").Value = 795.4545000000001
$ws.Range("M27").Value = -2026.3333
$ws.Range("N27").Value = -1009.4545

# Row 46
$ws.Range("H46").Value = 1681
$ws.Range("I46").Value = 970.2857
$ws.Range("J46").Value = 2676
$ws.Range("K46").Value = 970.2857
$ws.Range("L46").Value = 2676
$ws.Range("M46").Value = -782.2857
$ws.Range("N46").Value = -3052

# Row 61
$ws.Range("H61").Value = 1787.0769
$ws.Range("I61").Value = 953.5
$ws.Range("J61").Value = 3120.8
$ws.Range("K61").Value = 953.5
$ws.Range("L61").Value = 3120.8
$ws.Range("M61").Value = -751.5
$ws.Range("N61").Value = -3524.8

# Row 93
$ws.Range("H93").Value = 1277.6154
$ws.Range("I93").Value = 1272.7142
$ws.Range("J93").Value = 1283.3334
$ws.Range("K93").Value = 1272.7142
$ws.Range("L93").Value = 1283.3334
$ws.Range("M93").Value = -24.71419999999989
$ws.Range("N93").Value = -3779.3334

# Row 113
$ws.Range("H113").Value = 1787.0769
$ws.Range("I113").Value = 953.5
$ws.Range("J113").Value = 3120.8
$ws.Range("K113").Value = 953.5
$ws.Range("L113").Value = 3120.8
$ws.Range("M113").Value = 1216.5
$ws.Range("N113").Value = -7460.8

# Row 122
$ws.Range("H122").Value = 6790.2856
$ws.Range("I122").Value = 9392.076999999999
$ws.Range("J122").Value = 2562.375
$ws.Range("K122").Value = 28176.231
$ws.Range("L122").Value = 7687.125
$ws.Range("M122").Value = -25726.231
$ws.Range("N122").Value = -12587.125

$ws = $wb.Worksheets.Item("WVR")
# Row 39
$ws.Range("H39").Value = 59049
$ws.Range("J39").Value = 59049
$ws.Range("L39").Value = 59049
$ws.Range("N39").Value = -59875

# Row 81
$ws.Range("H81").Value = 1276.6154
$ws.Range("I81").Value = 1236
$ws.Range("J81").Value = 1324
$ws.Range("K81").Value = 2472
$ws.Range("L81").Value = 2648
$ws.Range("M81").Value = -1411
$ws.Range("N81").Value = -4770

# Row 84
$ws.Range("H84").Value = 1276.6154
$ws.Range("I84").Value = 1236
$ws.Range("J84").Value = 1324
$ws.Range("K84").Value = 12360
$ws.Range("L84").Value = 13240
$ws.Range("M84").Value = -7056
$ws.Range("N84").Value = -23848

# Row 136
$ws.Range("H136").Value = 25003322
$ws.Range("I136").Value = 55557384
$ws.Range("J136").Value = 4545.4546
$ws.Range("K136").Value = 166672152
$ws.Range("L136").Value = 13636.3638
$ws.Range("M136").Value = -166669602
$ws.Range("N136").Value = -18736.3638
